# Update crypto price and volume(1h) values for rows 2-51 (columns D and E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "0.998", "305.60") are preserved as text, not converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.933.08"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.534.21"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "305.60"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "102.14"
$ws.Range("E6").Value = "  +8.02%  "
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "37.95"
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "2.921.82"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "2.504.85"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "15.19"
$ws.Range("E16").Value = "  +6.98%  "
$ws.Range("D17").Value = "0.872"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "42.945.66"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").Value = "  +4.06%  "
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "6.51"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "71.63"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "252.88"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "2.06"
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("D26").Value = "27.20"
$ws.Range("E26").Value = "  -6.06%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +8.72%  "
$ws.Range("D29").Value = "10.33"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "39.22"
$ws.Range("E30").Value = "  +5.51%  "
$ws.Range("D31").Value = "6.18"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").Value = "157.69"
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "0.0797"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  -3.95%  "
$ws.Range("D36").Value = "3.27"
$ws.Range("E36").Value = "  -3.46%  "
$ws.Range("D37").Value = "18.41"
$ws.Range("E37").Value = "  +4.84%  "
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("D39").Value = "24.16"
$ws.Range("E39").Value = "  +5.28%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("D42").Value = "2.12"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("D43").Value = "3.90"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").Value = "0.0305"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "2.046.10"
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").Value = "86.44"
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").Value = "8.96"
$ws.Range("E48").Value = "  -3.84%  "
$ws.Range("D49").Value = "2.782.75"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").Value = "0.193"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("D51").Value = "103.00"
$ws.Range("E51").Value = "  -2.73%  "
